$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: E2 h5 -> "No es posible ingresar"
$ws.Range("E2").Value = "No es posible ingresar"

# Row 3: C3 div -> long "no results" message
$ws.Range("C3").Value = "No encontramos resultados para tu consulta. Te recomendamos usar frases cortas y palabras claves. Ej.: ""caja de ahorro""."

# Row 4: clear C4 and E4 (Dato2CP3 / Dato3CP3 removed)
$ws.Range("C4").ClearContents()
$ws.Range("E4").ClearContents()

# Row 5: new CP004 case data
$ws.Range("A5").Value = "CP004_Sucursal006"
$ws.Range("B5").Value = "FLORES"
$ws.Range("C5").Value = 6
$ws.Range("D5").Value = "RIVADAVIA"
$ws.Range("E5").Value = "FLORES"

# Rows 6-9: keep CP005..CP008 labels in column A, clear B/C/E data columns
$ws.Range("A6").Value = "CP005"
$ws.Range("B6").ClearContents()
$ws.Range("C6").ClearContents()
$ws.Range("E6").ClearContents()

$ws.Range("A7").Value = "CP006"
$ws.Range("B7").ClearContents()
$ws.Range("C7").ClearContents()
$ws.Range("E7").ClearContents()

$ws.Range("A8").Value = "CP007"
$ws.Range("B8").ClearContents()
$ws.Range("C8").ClearContents()
$ws.Range("E8").ClearContents()

$ws.Range("A9").Value = "CP008"
$ws.Range("B9").ClearContents()
$ws.Range("C9").ClearContents()
$ws.Range("E9").ClearContents()

# Column D width change (11.140625 -> 20.42578125)
$ws.Columns.Item(4).ColumnWidth = 20.42578125

# Update selection to F6 (matches recorded cursor position after edits)
$ws.Range("F6").Select()
